# Daily attendance processing - 2025-12-03 17:03:27
# Reverses the order of the comma-separated "Recorded By" entries in column G.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ", "
        $n = $parts.Count
        if ($n -gt 1) {
            $rev = @()
            for ($i = $n - 1; $i -ge 0; $i--) {
                $rev += $parts[$i]
            }
            $cell.Value2 = [string]::Join(", ", $rev)
        }
    }
}
